$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.791.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.582.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.045.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.656.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.591.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.706.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0811"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "456.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "157.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  -0.96%  "
